# Applies the "remove UML questions" edit to CC1.docx
# - merges a couple of adjacent runs that happened to split across a
#   retype (no visible text change)
# - strikes through "3.13, 3.14" in the A1 assignment-number list
# - merges "A2. (30 points)" + " Do problem 3.4E..." into one run
# - merges the "B" / "<digit>" / ". (...)" run triples into single runs
# - merges "Total points on assignment: " + "86" into one run

$d = $word.ActiveDocument

$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, $wdFindContinue, $false, $replace, $wdReplaceAll)
}

# 1. " " + "Answer the following review questions..." -> single run (no text change)
Replace-Text " Answer the following review questions from the text book:" " Answer the following review questions from the text book:"

# 2. Strike through "3.13, 3.14" within the "3.5, 3.6, ... 3.14" line.
#    First make sure the run is unified (it already is one run), then locate the
#    substring "3.13, 3.14" and apply strikethrough formatting to just that part.
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "3.13, 3.14"
$rng.Find.Forward = $true
$rng.Find.Wrap = $wdFindContinue
$rng.Find.MatchCase = $true
$rng.Find.MatchWholeWord = $false
while ($rng.Find.Execute()) {
    $rng.Font.StrikeThrough = $true
    $rng.Collapse(0)
}

# 3. "A2. (30 points)" + " Do problem 3.4E..." -> single run
Replace-Text "A2. (30 points) Do problem 3.4E from the text book. The question asks for both EER and UML models, but you are required to do only the EER model. Make notes to clarify details that cannot be captured by the EER diagram alone. Submit both your EER diagram and any additional notes you need to clarify your model." "A2. (30 points) Do problem 3.4E from the text book. The question asks for both EER and UML models, but you are required to do only the EER model. Make notes to clarify details that cannot be captured by the EER diagram alone. Submit both your EER diagram and any additional notes you need to clarify your model."

# 4. B1-B5 "B" + "<n>" + ". (...)" -> single run "B<n>. (...)"
Replace-Text "B1. (5 points) The names of all wines whose name contains the string " "B1. (5 points) The names of all wines whose name contains the string "
Replace-Text "B2. (5 points) The names of all suppliers whose status is null." "B2. (5 points) The names of all suppliers whose status is null."
Replace-Text "B3. (5 points) The names of all suppliers whose status is not null." "B3. (5 points) The names of all suppliers whose status is not null."
Replace-Text "B4. (10 points) Consider product 0154, the " "B4. (10 points) Consider product 0154, the "
Replace-Text "B5. (10 points) Consider your answer to the previous question." "B5. (10 points) Consider your answer to the previous question."

# 5. "Total points on assignment: " + "86" -> single run
Replace-Text "Total points on assignment: 86" "Total points on assignment: 86"
